$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 832.6667
$ws.Range("I32").Value = 501
$ws.Range("J32").Value = 998.5
$ws.Range("K32").Value = 501
$ws.Range("L32").Value = 998.5
$ws.Range("M32").Value = -175
$ws.Range("N32").Value = -1650.5
$ws.Range("H98").Value = 3821.6956
$ws.Range("I98").Value = 1391.5834
$ws.Range("J98").Value = 6472.727
$ws.Range("K98").Value = 1391.5834
$ws.Range("L98").Value = 6472.727
$ws.Range("M98").Value = 106.4166
$ws.Range("N98").Value = -9468.726999999999
$ws.Range("H108").Value = 39060
$ws.Range("J108").Value = 39060
$ws.Range("L108").Value = 39060
$ws.Range("N108").Value = -46740
$ws.Range("H122").Value = 3821.6956
$ws.Range("I122").Value = 1391.5834
$ws.Range("J122").Value = 6472.727
$ws.Range("K122").Value = 4174.7502
$ws.Range("L122").Value = 19418.181
$ws.Range("M122").Value = -1724.7502
$ws.Range("N122").Value = -24318.181
$ws.Range("H129").Value = 842.05
$ws.Range("J129").Value = 872.4894
$ws.Range("L129").Value = 2617.4682
$ws.Range("N129").Value = -12617.4682
$ws.Range("H137").Value = 2072578
$ws.Range("I137").Value = 2977583.8
$ws.Range("J137").Value = 3993.5715
$ws.Range("K137").Value = 8932751.399999999
$ws.Range("L137").Value = 11980.7145
$ws.Range("M137").Value = -8930201.399999999
$ws.Range("N137").Value = -17080.7145
$ws.Range("H141").Value = 17729.77
$ws.Range("I141").Value = 20044.273
$ws.Range("K141").Value = 60132.819
$ws.Range("M141").Value = -54952.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1827.1666
$ws.Range("I61").Value = 1323.5
$ws.Range("J61").Value = 2010.3182
$ws.Range("K61").Value = 1323.5
$ws.Range("L61").Value = 2010.3182
$ws.Range("M61").Value = -1111.5
$ws.Range("N61").Value = -2434.3182
$ws.Range("H74").Value = 6724.5293
$ws.Range("I74").Value = 8193.182000000001
$ws.Range("J74").Value = 4032
$ws.Range("K74").Value = 8193.182000000001
$ws.Range("L74").Value = 4032
$ws.Range("M74").Value = -7319.182000000001
$ws.Range("N74").Value = -5780
$ws.Range("H77").Value = 6724.5293
$ws.Range("I77").Value = 8193.182000000001
$ws.Range("J77").Value = 4032
$ws.Range("K77").Value = 40965.91
$ws.Range("L77").Value = 20160
$ws.Range("M77").Value = -36597.91
$ws.Range("N77").Value = -28896
$ws.Range("H132").Value = 2516.6843
$ws.Range("I132").Value = 1501.2858
$ws.Range("J132").Value = 5359.8
$ws.Range("K132").Value = 4503.857400000001
$ws.Range("L132").Value = 16079.4
$ws.Range("M132").Value = -1973.857400000001
$ws.Range("N132").Value = -21139.4
$ws.Range("H136").Value = 1827.1666
$ws.Range("I136").Value = 1323.5
$ws.Range("J136").Value = 2010.3182
$ws.Range("K136").Value = 3970.5
$ws.Range("L136").Value = 6030.9546
$ws.Range("M136").Value = -1420.5
$ws.Range("N136").Value = -11130.9546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20630
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22184
$ws.Range("H134").Value = 2386.375
$ws.Range("I134").Value = 1770.6897
$ws.Range("J134").Value = 8338
$ws.Range("K134").Value = 5312.0691
$ws.Range("L134").Value = 25014
$ws.Range("M134").Value = -2777.0691
$ws.Range("N134").Value = -30084

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2541.2559
$ws.Range("I31").Value = 1042.2778
$ws.Range("J31").Value = 3620.52
$ws.Range("K31").Value = 1042.2778
$ws.Range("L31").Value = 3620.52
$ws.Range("M31").Value = -747.2778000000001
$ws.Range("N31").Value = -4210.52
$ws.Range("H34").Value = 2541.2559
$ws.Range("I34").Value = 1042.2778
$ws.Range("J34").Value = 3620.52
$ws.Range("K34").Value = 1042.2778
$ws.Range("L34").Value = 3620.52
$ws.Range("M34").Value = -840.2778000000001
$ws.Range("N34").Value = -4024.52
$ws.Range("H58").Value = 2777.5715
$ws.Range("I58").Value = 1657.7322
$ws.Range("J58").Value = 7256.9287
$ws.Range("K58").Value = 1657.7322
$ws.Range("L58").Value = 7256.9287
$ws.Range("M58").Value = -1454.7322
$ws.Range("N58").Value = -7662.9287
$ws.Range("H74").Value = 33701
$ws.Range("J74").Value = 33701
$ws.Range("L74").Value = 33701
$ws.Range("N74").Value = -35449
$ws.Range("H77").Value = 33701
$ws.Range("J77").Value = 33701
$ws.Range("L77").Value = 101103
$ws.Range("N77").Value = -109839
$ws.Range("H125").Value = 35325
$ws.Range("J125").Value = 35325
$ws.Range("L125").Value = 35325
$ws.Range("N125").Value = -40245
$ws.Range("H132").Value = 2776.3462
$ws.Range("I132").Value = 1626.7778
$ws.Range("J132").Value = 5362.875
$ws.Range("K132").Value = 4880.3334
$ws.Range("L132").Value = 16088.625
$ws.Range("M132").Value = -2350.3334
$ws.Range("N132").Value = -21148.625
$ws.Range("H136").Value = 2777.5715
$ws.Range("I136").Value = 1657.7322
$ws.Range("J136").Value = 7256.9287
$ws.Range("K136").Value = 4973.196599999999
$ws.Range("L136").Value = 21770.7861
$ws.Range("M136").Value = -2423.196599999999
$ws.Range("N136").Value = -26870.7861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 818.8
$ws.Range("I44").Value = 383.2857
$ws.Range("J44").Value = 1199.875
$ws.Range("K44").Value = 1149.8571
$ws.Range("L44").Value = 3599.625
$ws.Range("M44").Value = -751.8571000000002
$ws.Range("N44").Value = -4395.625
$ws.Range("H113").Value = 4167252.2
$ws.Range("J113").Value = 8333885
$ws.Range("L113").Value = 25001655
$ws.Range("N113").Value = -25005995
$ws.Range("H131").Value = 679.4
$ws.Range("I131").Value = 254.65218
$ws.Range("J131").Value = 806.2727
$ws.Range("K131").Value = 763.9565399999999
$ws.Range("L131").Value = 2418.8181
$ws.Range("M131").Value = 4276.04346
$ws.Range("N131").Value = -12498.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51622
$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -158112
$ws.Range("H80").Value = 50003580
$ws.Range("J80").Value = 4950
$ws.Range("L80").Value = 4950
$ws.Range("N80").Value = -6946
$ws.Range("H83").Value = 50003580
$ws.Range("J83").Value = 4950
$ws.Range("L83").Value = 24750
$ws.Range("N83").Value = -34734
$ws.Range("H132").Value = 6349
$ws.Range("I132").Value = 5248
$ws.Range("J132").Value = 7083
$ws.Range("K132").Value = 15744
$ws.Range("L132").Value = 21249
$ws.Range("M132").Value = -13214
$ws.Range("N132").Value = -26309

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1747.6
$ws.Range("I22").Value = 1486.3636
$ws.Range("J22").Value = 2066.889
$ws.Range("K22").Value = 1486.3636
$ws.Range("L22").Value = 2066.889
$ws.Range("M22").Value = -1191.3636
$ws.Range("N22").Value = -2656.889
$ws.Range("H27").Value = 1747.6
$ws.Range("I27").Value = 1486.3636
$ws.Range("J27").Value = 2066.889
$ws.Range("K27").Value = 1486.3636
$ws.Range("L27").Value = 2066.889
$ws.Range("M27").Value = -1379.3636
$ws.Range("N27").Value = -2280.889
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H68").Value = 692.30304
$ws.Range("I68").Value = 692.30304
$ws.Range("K68").Value = 692.30304
$ws.Range("M68").Value = 56.69695999999999
$ws.Range("H71").Value = 692.30304
$ws.Range("I71").Value = 692.30304
$ws.Range("K71").Value = 3461.5152
$ws.Range("M71").Value = 282.4848000000002
$ws.Range("H82").Value = 4051.1714
$ws.Range("I82").Value = 5291.909
$ws.Range("J82").Value = 1951.4615
$ws.Range("K82").Value = 5291.909
$ws.Range("L82").Value = 1951.4615
$ws.Range("M82").Value = -4930.909
$ws.Range("N82").Value = -2673.4615
$ws.Range("H85").Value = 4051.1714
$ws.Range("I85").Value = 5291.909
$ws.Range("J85").Value = 1951.4615
$ws.Range("K85").Value = 5291.909
$ws.Range("L85").Value = 1951.4615
$ws.Range("M85").Value = -4043.909
$ws.Range("N85").Value = -4447.461499999999
$ws.Range("H112").Value = 32500
$ws.Range("J112").Value = 32500
$ws.Range("L112").Value = 32500
$ws.Range("N112").Value = -35454
$ws.Range("H122").Value = 9282.923000000001
$ws.Range("J122").Value = 9971.429
$ws.Range("L122").Value = 29914.287
$ws.Range("N122").Value = -34814.287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15687.375
$ws.Range("I62").Value = 2375
$ws.Range("J62").Value = 28999.75
$ws.Range("K62").Value = 2375
$ws.Range("L62").Value = 28999.75
$ws.Range("M62").Value = -1751
$ws.Range("N62").Value = -30247.75
$ws.Range("H65").Value = 15687.375
$ws.Range("I65").Value = 2375
$ws.Range("J65").Value = 28999.75
$ws.Range("K65").Value = 11875
$ws.Range("L65").Value = 144998.75
$ws.Range("M65").Value = -8755
$ws.Range("N65").Value = -151238.75
$ws.Range("H82").Value = 31333.334
$ws.Range("J82").Value = 42500
$ws.Range("L82").Value = 42500
$ws.Range("N82").Value = -43266
$ws.Range("H85").Value = 31333.334
$ws.Range("J85").Value = 42500
$ws.Range("L85").Value = 42500
$ws.Range("N85").Value = -45152
$ws.Range("H122").Value = 4163.5713
$ws.Range("I122").Value = 1988.5
$ws.Range("K122").Value = 5965.5
$ws.Range("M122").Value = -3515.5
$ws.Range("H136").Value = 12065.0625
$ws.Range("I136").Value = 12495.889
$ws.Range("K136").Value = 37487.667
$ws.Range("M136").Value = -34937.667
